$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update deaths for the week of 2020-03-28 (row 3)
$ws.Range("B3").Value = 87000

# Row 5 (week of 2020-03-14): derive inferred cases from confirmed cases x 100,
# and back out the implied lag multiplier instead of hard-coding it.
$ws.Range("C5").Formula = "=100*B3"
$ws.Range("J5").Formula = "=C5/I5"
